# Shrink the four axis-label textboxes ("EBVhost", "Tmax", "DM (pc)", "dm15")
# on the covariance-plot slide: font size 16pt -> 12pt, and re-fit their
# boxes to the smaller text (PowerPoint's auto-fit shrinks the box height;
# rotated boxes shrink about their centre, the "dm15" box also shifts right
# a touch, matching the hand-nudged layout from the commit).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Target geometry (EMU, taken from the rendered layout) per textbox name.
$targets = @{
    "TextBox 4" = @{ Left = 296.8447244094488;  Top = 217.82700787401575; Width = 123.9467736535433; Height = 21.810944881889764 }  # EBVhost
    "TextBox 6" = @{ Left = 296.8447244094488;  Top = 119.25464566929134; Width = 123.9467736535433; Height = 21.810944881889764 }  # Tmax
    "TextBox 7" = @{ Left = 397.09449818897633; Top = 302.57087714173224; Width = 123.9467736535433; Height = 21.810944881889764 }  # DM (pc)
    "TextBox 8" = @{ Left = 526.9056095511811;  Top = 302.5707874015748;  Width = 123.9467736535433; Height = 21.810944881889764 }  # dm15
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $name = $shp.Name
    if ($targets.ContainsKey($name)) {
        $t = $targets[$name]

        # Shrink the run text from 16pt to 12pt.
        $shp.TextFrame.TextRange.Font.Size = 12

        # Re-fit the shape's box to the new (smaller) auto-fit text size.
        $shp.Left = $t.Left
        $shp.Top = $t.Top
        $shp.Width = $t.Width
        $shp.Height = $t.Height
    }
}
